$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33-35 reorder: Coin name + Link columns ---
$ws.Range("B33").Value = 'TrustWalletToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

# --- Price (column D) updates: force Text format to preserve exact digits/trailing zeros ---
$dCells = "D2", "D3", "D5", "D8", "D12", "D13", "D16", "D17", "D31", "D33", "D34", "D35", "D39", "D40", "D42", "D44", "D45", "D47", "D49"
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value = '26.918.93'
$ws.Range("D3").Value = '1.551.00'
$ws.Range("D5").Value = '206.60'
$ws.Range("D8").Value = '22.01'
$ws.Range("D12").Value = '1.772.30'
$ws.Range("D13").Value = '1.545.85'
$ws.Range("D16").Value = '26.908.95'
$ws.Range("D17").Value = '61.58'
$ws.Range("D31").Value = '1.09'
$ws.Range("D33").Value = '1.15'
$ws.Range("D34").Value = '1.418.91'
$ws.Range("D35").Value = '3.10'
$ws.Range("D39").Value = '0.527'
$ws.Range("D40").Value = '0.805'
$ws.Range("D42").Value = '5.63'
$ws.Range("D44").Value = '0.996'
$ws.Range("D45").Value = '64.43'
$ws.Range("D47").Value = '1.686.21'
$ws.Range("D49").Value = '0.0520'

# --- Volume(1h) percentage text (column E) updates ---
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  -0.64%  '
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  -0.62%  '
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("E18").Value = '  +2.82%  '
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -1.29%  '
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("E29").Value = '  -0.58%  '
$ws.Range("E30").Value = '  +1.41%  '
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("E33").Value = '  +20.72%  '
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("E35").Value = '  +3.51%  '
$ws.Range("E36").Value = '  +2.93%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("E42").Value = '  +2.21%  '
$ws.Range("E43").Value = '  +1.69%  '
$ws.Range("E44").Value = '  +0.73%  '
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("E48").Value = '  +1.09%  '
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("E50").Value = '  +5.21%  '
$ws.Range("E51").Value = '  +0.25%  '
